$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per the crypto price/volume refresh.
# D-column values that look like plain numbers are forced to Text format
# first (matching the source data which stores them as strings), then the
# cell style is reset to Normal so no stray style index is left behind.

$ws.Range('D2').Value = '68.954.63'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '2.669.00'
$ws.Range('E3').Value = '  +3.04%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.23%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').Value = '2.669.22'
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('E10').Value = '  +14.45%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.351'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000189'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.99%  '
$ws.Range('D16').Value = '3.153.55'
$ws.Range('E16').Value = '  +2.98%  '
$ws.Range('D17').Value = '68.834.13'
$ws.Range('E17').Value = '  +2.75%  '
$ws.Range('D18').Value = '2.669.99'
$ws.Range('E18').Value = '  +2.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '368.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('E24').Value = '  +5.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000107'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.01%  '
$ws.Range('D29').Value = '2.808.79'
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '583.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  +5.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.87'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('E35').Value = '  +6.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '158.83'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('E39').Value = '  +6.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.43'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.63%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.370'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.70%  '
$ws.Range('E43').Value = '  +8.49%  '
$ws.Range('E44').Value = '  +5.92%  '
$ws.Range('E45').Value = '  +15.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.67%  '
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('E50').Value = '  +3.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.91%  '

Write-Output "Applied 82 cell updates to cryptos sheet."
